$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 10
$ws.Range("AE2").Value = 17
$ws.Range("AG2").Value = 26
$ws.Range("AI2").Value = 67
$ws.Range("K3").Value = 8.5
$ws.Range("AD4").Value = 301
$ws.Range("G5").Value = 1.27
$ws.Range("H5").Value = 6.25
$ws.Range("I5").Value = 9
$ws.Range("N5").Value = 1.5
$ws.Range("O5").Value = 2.5
$ws.Range("U5").Value = 7
$ws.Range("G6").Value = 2.1
$ws.Range("I6").Value = 3.5
$ws.Range("P6").Value = 1.33
$ws.Range("Q6").Value = 3.25
$ws.Range("U6").Value = 11
$ws.Range("W6").Value = 19
$ws.Range("X6").Value = 15
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 19
$ws.Range("AI6").Value = 26
$ws.Range("G7").Value = 1.65
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 5.25
$ws.Range("L7").Value = 1.18
$ws.Range("M7").Value = 4.5
$ws.Range("N7").Value = 1.65
$ws.Range("O7").Value = 2.2
$ws.Range("P7").Value = 1.3
$ws.Range("Q7").Value = 3.4
$ws.Range("X7").Value = 12
$ws.Range("AA7").Value = 7.5
$ws.Range("G9").Value = 3.2
$ws.Range("J9").Value = 1.11
$ws.Range("K9").Value = 6.5
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4
$ws.Range("N11").Value = 1.75
$ws.Range("O11").Value = 2.05
$ws.Range("P11").Value = 1.33
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.75
$ws.Range("S11").Value = 2
$ws.Range("T11").Value = 8
$ws.Range("AC11").Value = 41
$ws.Range("G12").Value = 2.7
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 2.55
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 13
$ws.Range("X12").Value = 23
$ws.Range("Z12").Value = 8.5
$ws.Range("AE12").Value = 8
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 10
$ws.Range("AH12").Value = 26
$ws.Range("I13").Value = 1.8
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9
$ws.Range("N13").Value = 2.25
$ws.Range("O13").Value = 1.62
$ws.Range("V13").Value = 17
$ws.Range("AD13").Value = 501
$ws.Range("AF13").Value = 7.5
$ws.Range("G14").Value = 2.8
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 2.25
$ws.Range("L14").Value = 1.17
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 1.6
$ws.Range("O14").Value = 2.3
$ws.Range("P14").Value = 1.3
$ws.Range("Q14").Value = 3.4
$ws.Range("R14").Value = 1.53
$ws.Range("S14").Value = 2.38
$ws.Range("T14").Value = 13
$ws.Range("Y14").Value = 23
$ws.Range("Z14").Value = 15
$ws.Range("AA14").Value = 7
$ws.Range("AC14").Value = 34
$ws.Range("AD14").Value = 101
$ws.Range("AE14").Value = 11
$ws.Range("AJ14").Value = 21
$ws.Range("G15").Value = 2.5
$ws.Range("I15").Value = 2.75
$ws.Range("U15").Value = 12
$ws.Range("W15").Value = 26
$ws.Range("AD15").Value = 301
$ws.Range("AE15").Value = 8
$ws.Range("AI15").Value = 23
$ws.Range("K17").Value = 9.5
$ws.Range("G20").Value = 1.33
$ws.Range("H20").Value = 4.75
$ws.Range("I20").Value = 8.5
$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.98
$ws.Range("O20").Value = 1.88
$ws.Range("R20").Value = 2.25
$ws.Range("S20").Value = 1.57
$ws.Range("T20").Value = 5.5
$ws.Range("U20").Value = 5.5
$ws.Range("V20").Value = 9
$ws.Range("W20").Value = 8
$ws.Range("Z20").Value = 9.5
$ws.Range("AA20").Value = 9
$ws.Range("AB20").Value = 26
$ws.Range("AE20").Value = 19
$ws.Range("AG20").Value = 26
$ws.Range("AH20").Value = 126
$ws.Range("AI20").Value = 67
$ws.Range("AJ20").Value = 67
$ws.Range("L21").Value = 1.22
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = 1.73
$ws.Range("O21").Value = 2.08
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 1.75
$ws.Range("U21").Value = 34
$ws.Range("Z21").Value = 12
$ws.Range("AB21").Value = 21
$ws.Range("AD21").Value = 401
$ws.Range("AF21").Value = 6.5
$ws.Range("G22").Value = 2.9
$ws.Range("I22").Value = 2.2
$ws.Range("R22").Value = 1.62
$ws.Range("S22").Value = 2.2
$ws.Range("W22").Value = 29
$ws.Range("X22").Value = 21
$ws.Range("Y22").Value = 26
$ws.Range("AE22").Value = 9.5
$ws.Range("AF22").Value = 12
$ws.Range("G23").Value = 2.75
$ws.Range("H23").Value = 3.5
$ws.Range("AA23").Value = 7
$ws.Range("G24").Value = 3.4
$ws.Range("N24").Value = 2.5
$ws.Range("O24").Value = 1.5
$ws.Range("P24").Value = 1.57
$ws.Range("Q24").Value = 2.25
$ws.Range("R24").Value = 2.1
$ws.Range("S24").Value = 1.67
$ws.Range("T24").Value = 7.5
$ws.Range("AB24").Value = 19
$ws.Range("AC24").Value = 67
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 3.5
$ws.Range("I25").Value = 3.7
$ws.Range("P25").Value = 1.4
$ws.Range("Q25").Value = 2.75
$ws.Range("R25").Value = 1.8
$ws.Range("S25").Value = 1.91
$ws.Range("T25").Value = 7.5
$ws.Range("U25").Value = 9.5
$ws.Range("V25").Value = 9
$ws.Range("Y25").Value = 26
$ws.Range("Z25").Value = 10
$ws.Range("AI25").Value = 29
$ws.Range("AJ25").Value = 34
$ws.Range("K27").Value = 17
$ws.Range("P27").Value = 1.29
$ws.Range("Q27").Value = 3.5
$ws.Range("T27").Value = 12
$ws.Range("V27").Value = 10
$ws.Range("Z27").Value = 17
$ws.Range("AE27").Value = 13
$ws.Range("G28").Value = 3.7
$ws.Range("H28").Value = 3.6
$ws.Range("I28").Value = 1.95
$ws.Range("J28").Value = 1.05
$ws.Range("K28").Value = 11
$ws.Range("N28").Value = 1.88
$ws.Range("O28").Value = 1.98
$ws.Range("T28").Value = 11
$ws.Range("U28").Value = 19
$ws.Range("Z28").Value = 11
$ws.Range("AJ28").Value = 26